# Apply crypto price / Volume(1h) updates scraped on Sat Apr 22 08:37:19 UTC 2023.
#
# Every target cell in this sheet holds its number/percentage as literal TEXT
# (inline strings such as "1.003" or "27.325.34", sometimes with multiple dots
# or padding spaces like "  -3.04%  "). A plain Range.Value assignment would let
# Excel auto-coerce numeric-looking text into a real Number, which would change
# the cell type and drop formatting such as the thousand-dot grouping or the
# padding spaces. To avoid that, each cell is briefly switched to the "Text"
# number format before the write and its original style is restored right after,
# so the only visible change is the cell content - exactly matching the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '27.325.34'
$cell.Style = $origStyle

$cell = $ws.Range("E2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.04%  '
$cell.Style = $origStyle

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.856.60'
$cell.Style = $origStyle

$cell = $ws.Range("E3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.84%  '
$cell.Style = $origStyle

$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.003'
$cell.Style = $origStyle

$cell = $ws.Range("E4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.29%  '
$cell.Style = $origStyle

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '324.27'
$cell.Style = $origStyle

$cell = $ws.Range("E5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.14%  '
$cell.Style = $origStyle

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.003'
$cell.Style = $origStyle

$cell = $ws.Range("E6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.16%  '
$cell.Style = $origStyle

$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.4526'
$cell.Style = $origStyle

$cell = $ws.Range("E7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.30%  '
$cell.Style = $origStyle

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.3869'
$cell.Style = $origStyle

$cell = $ws.Range("E8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.59%  '
$cell.Style = $origStyle

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '48.38'
$cell.Style = $origStyle

$cell = $ws.Range("E9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -8.66%  '
$cell.Style = $origStyle

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.07936'
$cell.Style = $origStyle

$cell = $ws.Range("E10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -6.12%  '
$cell.Style = $origStyle

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.015'
$cell.Style = $origStyle

$cell = $ws.Range("E11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.37%  '
$cell.Style = $origStyle

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '21.39'
$cell.Style = $origStyle

$cell = $ws.Range("E12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.95%  '
$cell.Style = $origStyle

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.846.47'
$cell.Style = $origStyle

$cell = $ws.Range("E13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -7.04%  '
$cell.Style = $origStyle

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.907'
$cell.Style = $origStyle

$cell = $ws.Range("E14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.49%  '
$cell.Style = $origStyle

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.125'
$cell.Style = $origStyle

$cell = $ws.Range("E15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -5.47%  '
$cell.Style = $origStyle

$cell = $ws.Range("E16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.49%  '
$cell.Style = $origStyle

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '85.86'
$cell.Style = $origStyle

$cell = $ws.Range("E17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.96%  '
$cell.Style = $origStyle

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.00001029'
$cell.Style = $origStyle

$cell = $ws.Range("E18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.52%  '
$cell.Style = $origStyle

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.06558'
$cell.Style = $origStyle

$cell = $ws.Range("E19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.50%  '
$cell.Style = $origStyle

$cell = $ws.Range("E20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -5.99%  '
$cell.Style = $origStyle

$cell = $ws.Range("E21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.12%  '
$cell.Style = $origStyle

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.538'
$cell.Style = $origStyle

$cell = $ws.Range("E22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.26%  '
$cell.Style = $origStyle

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '27.315.50'
$cell.Style = $origStyle

$cell = $ws.Range("E23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.15%  '
$cell.Style = $origStyle

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.94'
$cell.Style = $origStyle

$cell = $ws.Range("E24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.41%  '
$cell.Style = $origStyle

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.287'
$cell.Style = $origStyle

$cell = $ws.Range("E25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.11%  '
$cell.Style = $origStyle

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.076.49'
$cell.Style = $origStyle

$cell = $ws.Range("E26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -5.76%  '
$cell.Style = $origStyle

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '153.62'
$cell.Style = $origStyle

$cell = $ws.Range("E27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.59%  '
$cell.Style = $origStyle

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '19.91'
$cell.Style = $origStyle

$cell = $ws.Range("E28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.19%  '
$cell.Style = $origStyle

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.068'
$cell.Style = $origStyle

$cell = $ws.Range("E29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.42%  '
$cell.Style = $origStyle

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.458'
$cell.Style = $origStyle

$cell = $ws.Range("E30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -5.53%  '
$cell.Style = $origStyle

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '121.06'
$cell.Style = $origStyle

$cell = $ws.Range("E31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.23%  '
$cell.Style = $origStyle

$cell = $ws.Range("B32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Stellar'
$cell.Style = $origStyle

$cell = $ws.Range("C32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell.Style = $origStyle

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.09301'
$cell.Style = $origStyle

$cell = $ws.Range("E32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.38%  '
$cell.Style = $origStyle

$cell = $ws.Range("B33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'ImmutableX'
$cell.Style = $origStyle

$cell = $ws.Range("C33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell.Style = $origStyle

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.9369'
$cell.Style = $origStyle

$cell = $ws.Range("E33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.82%  '
$cell.Style = $origStyle

$cell = $ws.Range("B34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'ARBITRUM'
$cell.Style = $origStyle

$cell = $ws.Range("C34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell.Style = $origStyle

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.461'
$cell.Style = $origStyle

$cell = $ws.Range("E34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.85%  '
$cell.Style = $origStyle

$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.598'
$cell.Style = $origStyle

$cell = $ws.Range("E35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.13%  '
$cell.Style = $origStyle

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.272'
$cell.Style = $origStyle

$cell = $ws.Range("E36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -5.64%  '
$cell.Style = $origStyle

$cell = $ws.Range("E37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.30%  '
$cell.Style = $origStyle

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.02223'
$cell.Style = $origStyle

$cell = $ws.Range("E38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.21%  '
$cell.Style = $origStyle

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.05990'
$cell.Style = $origStyle

$cell = $ws.Range("E39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.08%  '
$cell.Style = $origStyle

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.113'
$cell.Style = $origStyle

$cell = $ws.Range("E40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -11.33%  '
$cell.Style = $origStyle

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = $origStyle

$cell = $ws.Range("E41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.17%  '
$cell.Style = $origStyle

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5912'
$cell.Style = $origStyle

$cell = $ws.Range("E42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.46%  '
$cell.Style = $origStyle

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.1891'
$cell.Style = $origStyle

$cell = $ws.Range("E43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.73%  '
$cell.Style = $origStyle

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.11'
$cell.Style = $origStyle

$cell = $ws.Range("E44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -8.89%  '
$cell.Style = $origStyle

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.272'
$cell.Style = $origStyle

$cell = $ws.Range("E45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.30%  '
$cell.Style = $origStyle

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5625'
$cell.Style = $origStyle

$cell = $ws.Range("E46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.51%  '
$cell.Style = $origStyle

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.98'
$cell.Style = $origStyle

$cell = $ws.Range("E47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -6.44%  '
$cell.Style = $origStyle

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.373'
$cell.Style = $origStyle

$cell = $ws.Range("E48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.86%  '
$cell.Style = $origStyle

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.917'
$cell.Style = $origStyle

$cell = $ws.Range("E49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -6.24%  '
$cell.Style = $origStyle

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.06756'
$cell.Style = $origStyle

$cell = $ws.Range("E50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.71%  '
$cell.Style = $origStyle

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '108.83'
$cell.Style = $origStyle

$cell = $ws.Range("E51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.97%  '
$cell.Style = $origStyle

